$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 235
$ws1.Range("F3").Value = 1435
$ws1.Range("F4").Value = 20179
$ws1.Range("F5").Value = 800
$ws1.Range("F6").Value = 316
$ws1.Range("F9").Value = 7601
$ws1.Range("F10").Value = 519
$ws1.Range("F12").Value = 272
$ws1.Range("F13").Value = 39
$ws1.Range("F14").Value = 160
$ws1.Range("F15").Value = 125
$ws1.Range("F16").Value = 13
$ws1.Range("F18").Value = 194
$ws1.Range("F19").Value = 1341
$ws1.Range("F20").Value = 429
$ws1.Range("F21").Value = 74
$ws1.Range("F22").Value = 680
$ws1.Range("F23").Value = 51
$ws1.Range("F24").Value = 72
$ws1.Range("F26").Value = 325
$ws1.Range("F27").Value = 1112
$ws1.Range("F30").Value = 185
$ws1.Range("F31").Value = 5220
$ws1.Range("F32").Value = 564
$ws1.Range("F33").Value = 71
$ws1.Range("F34").Value = 2862
$ws1.Range("F38").Value = 12664
$ws1.Range("F39").Value = 1338
$ws1.Range("F40").Value = 86
$ws1.Range("F41").Value = 28
$ws1.Range("F42").Value = 57
$ws1.Range("F43").Value = 273
$ws1.Range("F44").Value = 374
$ws1.Range("F45").Value = 4004
$ws1.Range("F46").Value = 321
$ws1.Range("F47").Value = 94

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 192

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 235
$ws4.Range("F3").Value = 1435
$ws4.Range("F4").Value = 20179
$ws4.Range("F5").Value = 800
$ws4.Range("F6").Value = 316
$ws4.Range("F9").Value = 7601
$ws4.Range("F10").Value = 519
$ws4.Range("F12").Value = 272
$ws4.Range("F13").Value = 39
$ws4.Range("F14").Value = 160
$ws4.Range("F15").Value = 125
$ws4.Range("F16").Value = 13
$ws4.Range("F18").Value = 194
$ws4.Range("F19").Value = 1341
$ws4.Range("F20").Value = 429
$ws4.Range("F21").Value = 74
$ws4.Range("F22").Value = 680
$ws4.Range("F23").Value = 51
$ws4.Range("F24").Value = 72
$ws4.Range("F26").Value = 325
$ws4.Range("F27").Value = 1112
$ws4.Range("F30").Value = 185
$ws4.Range("F31").Value = 192
$ws4.Range("F32").Value = 564
$ws4.Range("F34").Value = 71
$ws4.Range("F36").Value = 2862
$ws4.Range("F40").Value = 12664
$ws4.Range("F41").Value = 1338
$ws4.Range("F42").Value = 86
$ws4.Range("F43").Value = 28
$ws4.Range("F44").Value = 57
$ws4.Range("F45").Value = 273
$ws4.Range("F46").Value = 374
$ws4.Range("F47").Value = 4004
$ws4.Range("F48").Value = 321
$ws4.Range("F49").Value = 94

Write-Host "Update complete"
